# Add the "Nouvelles fonctionnalites du site" section to the end of the
# document, replacing the last (second) empty trailing paragraph with:
#   - a centered bold/underlined size-32 heading
#   - an empty bold/underlined paragraph
#   - a "Page saison et episode" paragraph
#   - a bulleted list item "Mettre un systeme de commentaire"

$d = $word.ActiveDocument

# The last paragraph in the body is the trailing empty paragraph right
# before the section break; its Range is where the new content goes.
$target = $d.Paragraphs.Item($d.Paragraphs.Count).Range

$newContentXml = @'
<?xml version="1.0" encoding="utf-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>Nouvelles fonctionnalités du site</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:r><w:t xml:space="preserve">Page saison et </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>episode</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Mettre un système de commentaire</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($newContentXml)
